# Work diary so far — append the new work-log rows (58-63) covering the
# 2022-05-17 (serial 44698) session, resize the table to match, and update
# the sheet view selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$workDate = 44698

# --- Dates, Type (Réalisation) and Durée for the new rows ------------------
# (these all reuse already-interned shared strings / are plain numbers, so
# the order they're written in does not create any new shared-string
# entries.)
foreach ($r in 58..63) {
    $ws.Range("A$r").Value = $workDate
}
foreach ($r in 58..62) {
    $ws.Range("B$r").Value = "Réalisation"
}
$ws.Range("C58").Value = 1
$ws.Range("C59").Value = 0.5
$ws.Range("C60").Value = 1
$ws.Range("C61").Value = 2
$ws.Range("C62").Value = 0.5

# --- Descriptions / remarks / commit id ------------------------------------
# Written in this precise order so new shared-string entries are interned in
# the same sequence the original author typed them in.
$ws.Range("D58").Value = "Création du formulaire d'ajout d'articles"
$ws.Range("D59").Value = "Récupération des catégories dans la base de donnée"
$ws.Range("D60").Value = "Ajout d'articles dans la base de donnée"
$ws.Range("D61").Value = "Ajout d'image dans la base de donnée"
$ws.Range("E61").Value = "Doit encore modifier la base de donnée et insérer un autreur"
$ws.Range("F61").Value = "0057e8865b5049d5d6bdcae276b3df47b63d9038"
$ws.Range("D62").Value = "Bugfix du choix de l'auteur"
$ws.Range("E62").Value = "Le formulaire d'ajout d'article n'aime pas les apostrophes ou autres caractères qui ne sont pas des chiffres et des lettres. Je ne fais pour l'instant aucun check de donnée"
$ws.Range("E58").Value = "Journée effectuée en Home Office"

# --- Formatting: copy the look of the last pre-existing row (57) onto the
# new date cells (column A uses a dd/mm/yyyy-style numeric format) and give
# the new commit-id cell (F61) the same look as the other commit-id cells. --
$ws.Range("A57").Copy()
$ws.Range("A58:A63").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("F49").Copy()
$ws.Range("F61").PasteSpecial(-4122) # xlPasteFormats (commit-id look)
$excel.CutCopyMode = $false

# Taller rows for the wrapped multi-line text, matching the authored sheet.
$ws.Rows.Item(61).RowHeight = 30
$ws.Rows.Item(62).RowHeight = 60

# --- Expand the table (ListObject) to cover the newly added rows ----------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:F63"))

# --- Update the view so the new last row is the active selection ----------
$ws.Range("E58").Select()
$excel.ActiveWindow.ScrollRow = 40
